# orders_detail.xlsx -> update export excel bill
# Rewrites the jx:each template placeholders (order -> product export context),
# refreshes the comment documenting the jx:each() tag, centers all the
# template cells vertically, and moves the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Placeholder text rewrites (merged "value" cells in the template table).
#    Only the cells whose template expression actually changed are touched;
#    everything else (labels, headers, unchanged placeholders) is left as-is.
# ---------------------------------------------------------------------------
$ws.Range("D4").Value  = '${item.manufactureName}'
$ws.Range("D5").Value  = '${item.orderCode}'
$ws.Range("D7").Value  = '${item.customerPhone}'
$ws.Range("D8").Value  = '${item.customerAddress}'

$ws.Range("B14").Value = '${productDto.productName}'
$ws.Range("D14").Value = '${productDto.price}'
$ws.Range("F14").Value = '${productDto.discount}'
$ws.Range("H14").Value = '${productDto.quantity}'
$ws.Range("J14").Value = '${productDto.totalAmount}'

$ws.Range("H17").Value = '${item.amountShipping}'
$ws.Range("H18").Value = '${item.amountTotal}'
$ws.Range("H24").Value = '${item.founder}'

# ---------------------------------------------------------------------------
# 2. Update the jx:each(...) documentation comment anchored on A14 so it
#    matches the new iteration variable / list name / var name.
# ---------------------------------------------------------------------------
$comment = $ws.Range("A14").Comment
$newCommentText = "Microsoft Office User:`njx:each(items=""item.exportProductDetailDtoList"", var=""productDto"", lastCell=""L15"")`n"
$comment.Text($newCommentText)

# ---------------------------------------------------------------------------
# 3. Vertically center every cell that carries content (labels + values),
#    matching the "vertical center" alignment added across the template.
# ---------------------------------------------------------------------------
$xlCenter = -4108
$ranges = @(
    "B2:K2",
    "B4:C8", "D4:K8",
    "B9:K9",
    "B11:K11",
    "B13:K13",
    "B14:K14",
    "B16:G16", "H16:K16",
    "B17:G17", "H17:K17",
    "B18:G18", "H18:K18",
    "C20:E20", "H20:J20",
    "C21:E23", "H21:J23",
    "C24:E24", "H24:J24"
)
foreach ($addr in $ranges) {
    $ws.Range($addr).VerticalAlignment = $xlCenter
}

# ---------------------------------------------------------------------------
# 4. Move the active selection (cosmetic, mirrors the saved cursor position).
# ---------------------------------------------------------------------------
$ws.Range("K23").Select()
